$d = $word.ActiveDocument

# 1) Shorten the first-page title: drop the trailing "Proposal".
$d.Content.Find.Execute(
    "Standard Nondisclosure Agreement Proposal", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Standard Nondisclosure Agreement", 2)

# 2) Merge the "I / propose ... certify / are exactly the same as " runs
#    into a single run by replacing the (identically formatted) combined
#    text that spans them - Word's Find/Replace coalesces the runs it
#    rewrites into one when the replacement text matches like-for-like.
$d.Content.Find.Execute(
    "I propose that we sign a standard nondisclosure agreement on the following terms, which I certify are exactly the same as ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I propose that we sign a standard nondisclosure agreement on the following terms, which I certify are exactly the same as ",
    2)
